$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-13 Thursday" "2025-03-14 Friday"

Replace-Text "256×5=1280" "346×8=2768"
Replace-Text "849×5=4245" "225×8=1800"
Replace-Text "703×6=4218" "874×3=2622"
Replace-Text "917×8=7336" "612×9=5508"
Replace-Text "423×3=1269" "640×9=5760"
Replace-Text "243×7=1701" "564×6=3384"
Replace-Text "401×3=1203" "201×5=1005"
Replace-Text "156×7=1092" "366×2=732"
Replace-Text "972×4=3888" "548×8=4384"
Replace-Text "136×7=952" "866×5=4330"
Replace-Text "673×7=4711" "527×2=1054"
Replace-Text "536×5=2680" "435×9=3915"
Replace-Text "497×6=2982" "261×3=783"
Replace-Text "874×9=7866" "607×8=4856"
Replace-Text "674×5=3370" "876×4=3504"
Replace-Text "669×9=6021" "112×7=784"
Replace-Text "461×8=3688" "476×7=3332"
Replace-Text "567×3=1701" "333×5=1665"
Replace-Text "787×5=3935" "838×7=5866"
Replace-Text "131×7=917" "132×2=264"
Replace-Text "887×5=4435" "868×6=5208"
Replace-Text "756×5=3780" "430×3=1290"
Replace-Text "276×8=2208" "590×8=4720"
Replace-Text "956×2=1912" "610×9=5490"
Replace-Text "520×3=1560" "908×3=2724"
